$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.448.66'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.376.96'
$ws.Range('E3').Value = '  +5.03%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.31'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.652'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '71.55'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +12.58%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.471'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.75%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0980'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '56.86'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.96%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '27.19'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '2.727.45'
$ws.Range('E13').Value = '  +4.80%  '
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '16.08'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.28'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '2.375.67'
$ws.Range('E18').Value = '  +4.79%  '
$ws.Range('D19').Value = '43.439.56'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('E21').Value = '  +2.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.63'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '251.01'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.70'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +11.09%  '
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.03'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.23'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '22.71'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '174.37'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  +5.09%  '
$ws.Range('E32').Value = '  -5.93%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.08'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.45'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +6.60%  '
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.92'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '18.71'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +8.37%  '
$ws.Range('E44').Value = '  +8.34%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '100.62'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.52'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.27%  '
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0957'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.446.81'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000206'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -7.94%  '
$ws.Range('D51').Value = '2.602.33'
$ws.Range('E51').Value = '  +5.09%  '
